# Update the "Datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'Datos actualizados a 2 de Abril de 2020 a las 15:20'

# Refresh the country/casos data table (A4:H209), which has been
# re-sorted by "Casos totales" (column B) descending and updated with
# the latest counts.
$data = New-Object 'object[,]' 206,8
$data[0,0] = 'Estados Unidos'
$data[0,1] = 215357
$data[0,2] = 354
$data[0,3] = 8878
$data[0,4] = 201366
$data[0,5] = 5005
$data[0,6] = 11
$data[0,7] = 5113
$data[1,0] = 'Italia'
$data[1,1] = 110574
$data[1,2] = 0
$data[1,3] = 16847
$data[1,4] = 80572
$data[1,5] = 4035
$data[1,6] = 0
$data[1,7] = 13155
$data[2,0] = 'España'
$data[2,1] = 110238
$data[2,2] = 6120
$data[2,3] = 26743
$data[2,4] = 73492
$data[2,5] = 6092
$data[2,6] = 616
$data[2,7] = 10003
$data[3,0] = 'China'
$data[3,1] = 81589
$data[3,2] = 35
$data[3,3] = 76408
$data[3,4] = 1863
$data[3,5] = 429
$data[3,6] = 6
$data[3,7] = 3318
$data[4,0] = 'Alemania'
$data[4,1] = 79465
$data[4,2] = 1484
$data[4,3] = 19175
$data[4,4] = 59331
$data[4,5] = 3408
$data[4,6] = 28
$data[4,7] = 959
$data[5,0] = 'Francia'
$data[5,1] = 56989
$data[5,2] = 0
$data[5,3] = 10935
$data[5,4] = 42022
$data[5,5] = 6017
$data[5,6] = 0
$data[5,7] = 4032
$data[6,0] = 'Iran'
$data[6,1] = 50468
$data[6,2] = 2875
$data[6,3] = 16711
$data[6,4] = 30597
$data[6,5] = 3956
$data[6,6] = 124
$data[6,7] = 3160
$data[7,0] = 'Reino Unido'
$data[7,1] = 33718
$data[7,2] = 4244
$data[7,3] = 135
$data[7,4] = 30662
$data[7,5] = 163
$data[7,6] = 569
$data[7,7] = 2921
$data[8,0] = 'Suiza'
$data[8,1] = 18267
$data[8,2] = 499
$data[8,3] = 4013
$data[8,4] = 13749
$data[8,5] = 348
$data[8,6] = 17
$data[8,7] = 505
$data[9,0] = 'Turquia'
$data[9,1] = 15679
$data[9,2] = 0
$data[9,3] = 333
$data[9,4] = 15069
$data[9,5] = 979
$data[9,6] = 0
$data[9,7] = 277
$data[10,0] = 'Belgica'
$data[10,1] = 15348
$data[10,2] = 1384
$data[10,3] = 2495
$data[10,4] = 11842
$data[10,5] = 1144
$data[10,6] = 183
$data[10,7] = 1011
$data[11,0] = 'Paises Bajos'
$data[11,1] = 14697
$data[11,2] = 1083
$data[11,3] = 250
$data[11,4] = 13108
$data[11,5] = 1053
$data[11,6] = 166
$data[11,7] = 1339
$data[12,0] = 'Austria'
$data[12,1] = 10927
$data[12,2] = 216
$data[12,3] = 1749
$data[12,4] = 9020
$data[12,5] = 227
$data[12,6] = 12
$data[12,7] = 158
$data[13,0] = 'Corea del Sur'
$data[13,1] = 9976
$data[13,2] = 89
$data[13,3] = 5828
$data[13,4] = 3979
$data[13,5] = 55
$data[13,6] = 4
$data[13,7] = 169
$data[14,0] = 'Canada'
$data[14,1] = 9731
$data[14,2] = 0
$data[14,3] = 1736
$data[14,4] = 7866
$data[14,5] = 120
$data[14,6] = 15
$data[14,7] = 129
$data[15,0] = 'Portugal'
$data[15,1] = 9034
$data[15,2] = 783
$data[15,3] = 68
$data[15,4] = 8757
$data[15,5] = 230
$data[15,6] = 22
$data[15,7] = 209
$data[16,0] = 'Brasil'
$data[16,1] = 6932
$data[16,2] = 52
$data[16,3] = 127
$data[16,4] = 6558
$data[16,5] = 296
$data[16,6] = 5
$data[16,7] = 247
$data[17,0] = 'Israel'
$data[17,1] = 6360
$data[17,2] = 268
$data[17,3] = 289
$data[17,4] = 6038
$data[17,5] = 107
$data[17,6] = 7
$data[17,7] = 33
$data[18,0] = 'Suecia'
$data[18,1] = 5466
$data[18,2] = 519
$data[18,3] = 103
$data[18,4] = 5081
$data[18,5] = 429
$data[18,6] = 43
$data[18,7] = 282
$data[19,0] = 'Australia'
$data[19,1] = 5137
$data[19,2] = 89
$data[19,3] = 345
$data[19,4] = 4767
$data[19,5] = 50
$data[19,6] = 2
$data[19,7] = 25
$data[20,0] = 'Noruega'
$data[20,1] = 5083
$data[20,2] = 206
$data[20,3] = 13
$data[20,4] = 5023
$data[20,5] = 105
$data[20,6] = 3
$data[20,7] = 47
$data[21,0] = 'Chequia'
$data[21,1] = 3604
$data[21,2] = 15
$data[21,3] = 61
$data[21,4] = 3503
$data[21,5] = 72
$data[21,6] = 1
$data[21,7] = 40
$data[22,0] = 'Rusia'
$data[22,1] = 3548
$data[22,2] = 771
$data[22,3] = 235
$data[22,4] = 3283
$data[22,5] = 8
$data[22,6] = 6
$data[22,7] = 30
$data[23,0] = 'Irlanda'
$data[23,1] = 3447
$data[23,2] = 0
$data[23,3] = 5
$data[23,4] = 3357
$data[23,5] = 103
$data[23,6] = 0
$data[23,7] = 85
$data[24,0] = 'Dinamarca'
$data[24,1] = 3355
$data[24,2] = 248
$data[24,3] = 1089
$data[24,4] = 2143
$data[24,5] = 153
$data[24,6] = 19
$data[24,7] = 123
$data[25,0] = 'Malasia'
$data[25,1] = 3116
$data[25,2] = 208
$data[25,3] = 767
$data[25,4] = 2299
$data[25,5] = 105
$data[25,6] = 5
$data[25,7] = 50
$data[26,0] = 'Chile'
$data[26,1] = 3031
$data[26,2] = 0
$data[26,3] = 234
$data[26,4] = 2781
$data[26,5] = 31
$data[26,6] = 0
$data[26,7] = 16
$data[27,0] = 'Ecuador'
$data[27,1] = 2758
$data[27,2] = 0
$data[27,3] = 58
$data[27,4] = 2602
$data[27,5] = 100
$data[27,6] = 0
$data[27,7] = 98
$data[28,0] = 'Rumania'
$data[28,1] = 2738
$data[28,2] = 278
$data[28,3] = 267
$data[28,4] = 2377
$data[28,5] = 78
$data[28,6] = 2
$data[28,7] = 94
$data[29,0] = 'Polonia'
$data[29,1] = 2633
$data[29,2] = 79
$data[29,3] = 56
$data[29,4] = 2532
$data[29,5] = 50
$data[29,6] = 2
$data[29,7] = 45
$data[30,0] = 'Filipinas'
$data[30,1] = 2633
$data[30,2] = 322
$data[30,3] = 51
$data[30,4] = 2475
$data[30,5] = 1
$data[30,6] = 11
$data[30,7] = 107
$data[31,0] = 'Japon'
$data[31,1] = 2384
$data[31,2] = 0
$data[31,3] = 472
$data[31,4] = 1855
$data[31,5] = 69
$data[31,6] = 0
$data[31,7] = 57
$data[32,0] = 'Luxemburgo'
$data[32,1] = 2319
$data[32,2] = 0
$data[32,3] = 80
$data[32,4] = 2210
$data[32,5] = 31
$data[32,6] = 0
$data[32,7] = 29
$data[33,0] = 'Pakistan'
$data[33,1] = 2291
$data[33,2] = 173
$data[33,3] = 107
$data[33,4] = 2153
$data[33,5] = 12
$data[33,6] = 4
$data[33,7] = 31
$data[34,0] = 'India'
$data[34,1] = 2032
$data[34,2] = 34
$data[34,3] = 150
$data[34,4] = 1824
$data[34,5] = 0
$data[34,6] = 0
$data[34,7] = 58
$data[35,0] = 'Arabia Saudita'
$data[35,1] = 1885
$data[35,2] = 165
$data[35,3] = 328
$data[35,4] = 1536
$data[35,5] = 31
$data[35,6] = 5
$data[35,7] = 21
$data[36,0] = 'Tailandia'
$data[36,1] = 1875
$data[36,2] = 104
$data[36,3] = 505
$data[36,4] = 1355
$data[36,5] = 23
$data[36,6] = 3
$data[36,7] = 15
$data[37,0] = 'Indonesia'
$data[37,1] = 1790
$data[37,2] = 113
$data[37,3] = 112
$data[37,4] = 1508
$data[37,5] = 0
$data[37,6] = 13
$data[37,7] = 170
$data[38,0] = 'Finlandia'
$data[38,1] = 1518
$data[38,2] = 72
$data[38,3] = 300
$data[38,4] = 1199
$data[38,5] = 62
$data[38,6] = 2
$data[38,7] = 19
$data[39,0] = 'Grecia'
$data[39,1] = 1415
$data[39,2] = 0
$data[39,3] = 52
$data[39,4] = 1311
$data[39,5] = 90
$data[39,6] = 1
$data[39,7] = 52
$data[40,0] = 'Sudafrica'
$data[40,1] = 1380
$data[40,2] = 0
$data[40,3] = 50
$data[40,4] = 1325
$data[40,5] = 7
$data[40,6] = 0
$data[40,7] = 5
$data[41,0] = 'Mexico'
$data[41,1] = 1378
$data[41,2] = 163
$data[41,3] = 35
$data[41,4] = 1306
$data[41,5] = 1
$data[41,6] = 8
$data[41,7] = 37
$data[42,0] = 'Peru'
$data[42,1] = 1323
$data[42,2] = 0
$data[42,3] = 394
$data[42,4] = 882
$data[42,5] = 49
$data[42,6] = 9
$data[42,7] = 47
$data[43,0] = 'Panama'
$data[43,1] = 1317
$data[43,2] = 0
$data[43,3] = 9
$data[43,4] = 1276
$data[43,5] = 50
$data[43,6] = 0
$data[43,7] = 32
$data[44,0] = 'Republica Dominicana'
$data[44,1] = 1284
$data[44,2] = 0
$data[44,3] = 9
$data[44,4] = 1218
$data[44,5] = 147
$data[44,6] = 0
$data[44,7] = 57
$data[45,0] = 'Islandia'
$data[45,1] = 1220
$data[45,2] = 0
$data[45,3] = 236
$data[45,4] = 982
$data[45,5] = 12
$data[45,6] = 0
$data[45,7] = 2
$data[46,0] = 'Serbia'
$data[46,1] = 1171
$data[46,2] = 111
$data[46,3] = 42
$data[46,4] = 1098
$data[46,5] = 62
$data[46,6] = 3
$data[46,7] = 31
$data[47,0] = 'Argentina'
$data[47,1] = 1133
$data[47,2] = 0
$data[47,3] = 248
$data[47,4] = 851
$data[47,5] = 0
$data[47,6] = 2
$data[47,7] = 34
$data[48,0] = 'Colombia'
$data[48,1] = 1065
$data[48,2] = 0
$data[48,3] = 39
$data[48,4] = 1009
$data[48,5] = 47
$data[48,6] = 0
$data[48,7] = 17
$data[49,0] = 'Croacia'
$data[49,1] = 1011
$data[49,2] = 48
$data[49,3] = 88
$data[49,4] = 916
$data[49,5] = 34
$data[49,6] = 1
$data[49,7] = 7
$data[50,0] = 'Singapur'
$data[50,1] = 1000
$data[50,2] = 0
$data[50,3] = 245
$data[50,4] = 751
$data[50,5] = 24
$data[50,6] = 1
$data[50,7] = 4
$data[51,0] = 'Eslovenia'
$data[51,1] = 897
$data[51,2] = 56
$data[51,3] = 70
$data[51,4] = 810
$data[51,5] = 31
$data[51,6] = 2
$data[51,7] = 17
$data[52,0] = 'Estonia'
$data[52,1] = 858
$data[52,2] = 79
$data[52,3] = 45
$data[52,4] = 802
$data[52,5] = 16
$data[52,6] = 6
$data[52,7] = 11
$data[53,0] = 'Argelia'
$data[53,1] = 847
$data[53,2] = 0
$data[53,3] = 61
$data[53,4] = 728
$data[53,5] = 0
$data[53,6] = 0
$data[53,7] = 58
$data[54,0] = 'Catar'
$data[54,1] = 835
$data[54,2] = 0
$data[54,3] = 71
$data[54,4] = 762
$data[54,5] = 37
$data[54,6] = 0
$data[54,7] = 2
$data[55,0] = 'Emiratos Arabes Unidos'
$data[55,1] = 814
$data[55,2] = 0
$data[55,3] = 61
$data[55,4] = 745
$data[55,5] = 2
$data[55,6] = 0
$data[55,7] = 8
$data[56,0] = 'Ucrania'
$data[56,1] = 804
$data[56,2] = 10
$data[56,3] = 13
$data[56,4] = 771
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 20
$data[57,0] = 'Hong Kong'
$data[57,1] = 802
$data[57,2] = 36
$data[57,3] = 154
$data[57,4] = 644
$data[57,5] = 8
$data[57,6] = 0
$data[57,7] = 4
$data[58,0] = 'Nueva Zelanda'
$data[58,1] = 797
$data[58,2] = 89
$data[58,3] = 92
$data[58,4] = 704
$data[58,5] = 2
$data[58,6] = 0
$data[58,7] = 1
$data[59,0] = 'Egipto'
$data[59,1] = 779
$data[59,2] = 0
$data[59,3] = 179
$data[59,4] = 548
$data[59,5] = 0
$data[59,6] = 0
$data[59,7] = 52
$data[60,0] = 'Irak'
$data[60,1] = 728
$data[60,2] = 0
$data[60,3] = 182
$data[60,4] = 494
$data[60,5] = 0
$data[60,6] = 0
$data[60,7] = 52
$data[61,0] = 'Crucero'
$data[61,1] = 712
$data[61,2] = 0
$data[61,3] = 619
$data[61,4] = 82
$data[61,5] = 15
$data[61,6] = 0
$data[61,7] = 11
$data[62,0] = 'Marruecos'
$data[62,1] = 676
$data[62,2] = 22
$data[62,3] = 29
$data[62,4] = 607
$data[62,5] = 1
$data[62,6] = 1
$data[62,7] = 40
$data[63,0] = 'Armenia'
$data[63,1] = 663
$data[63,2] = 92
$data[63,3] = 33
$data[63,4] = 625
$data[63,5] = 30
$data[63,6] = 1
$data[63,7] = 5
$data[64,0] = 'Lituania'
$data[64,1] = 649
$data[64,2] = 68
$data[64,3] = 7
$data[64,4] = 633
$data[64,5] = 11
$data[64,6] = 1
$data[64,7] = 9
$data[65,0] = 'Barein'
$data[65,1] = 635
$data[65,2] = 66
$data[65,3] = 341
$data[65,4] = 290
$data[65,5] = 3
$data[65,6] = 0
$data[65,7] = 4
$data[66,0] = 'Hungria'
$data[66,1] = 585
$data[66,2] = 60
$data[66,3] = 42
$data[66,4] = 522
$data[66,5] = 17
$data[66,6] = 1
$data[66,7] = 21
$data[67,0] = 'Bosnia y Herzegovina'
$data[67,1] = 518
$data[67,2] = 59
$data[67,3] = 20
$data[67,4] = 483
$data[67,5] = 4
$data[67,6] = 2
$data[67,7] = 15
$data[68,0] = 'Libano'
$data[68,1] = 494
$data[68,2] = 15
$data[68,3] = 43
$data[68,4] = 435
$data[68,5] = 3
$data[68,6] = 2
$data[68,7] = 16
$data[69,0] = 'Letonia'
$data[69,1] = 458
$data[69,2] = 12
$data[69,3] = 1
$data[69,4] = 457
$data[69,5] = 3
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = 'Bulgaria'
$data[70,1] = 449
$data[70,2] = 27
$data[70,3] = 25
$data[70,4] = 414
$data[70,5] = 17
$data[70,6] = 0
$data[70,7] = 10
$data[71,0] = 'Principado de Andorra'
$data[71,1] = 428
$data[71,2] = 38
$data[71,3] = 10
$data[71,4] = 403
$data[71,5] = 12
$data[71,6] = 1
$data[71,7] = 15
$data[72,0] = 'Eslovaquia'
$data[72,1] = 426
$data[72,2] = 26
$data[72,3] = 5
$data[72,4] = 420
$data[72,5] = 3
$data[72,6] = 0
$data[72,7] = 1
$data[73,0] = 'Tunez'
$data[73,1] = 423
$data[73,2] = 0
$data[73,3] = 5
$data[73,4] = 406
$data[73,5] = 10
$data[73,6] = 0
$data[73,7] = 12
$data[74,0] = 'Moldavia'
$data[74,1] = 423
$data[74,2] = 0
$data[74,3] = 23
$data[74,4] = 395
$data[74,5] = 65
$data[74,6] = 0
$data[74,7] = 5
$data[75,0] = 'Kazajistan'
$data[75,1] = 423
$data[75,2] = 43
$data[75,3] = 27
$data[75,4] = 393
$data[75,5] = 6
$data[75,6] = 0
$data[75,7] = 3
$data[76,0] = 'Azerbaiyan'
$data[76,1] = 400
$data[76,2] = 41
$data[76,3] = 26
$data[76,4] = 369
$data[76,5] = 7
$data[76,6] = 0
$data[76,7] = 5
$data[77,0] = 'Republica de Macedonia'
$data[77,1] = 384
$data[77,2] = 30
$data[77,3] = 17
$data[77,4] = 356
$data[77,5] = 4
$data[77,6] = 0
$data[77,7] = 11
$data[78,0] = 'Costa Rica'
$data[78,1] = 375
$data[78,2] = 0
$data[78,3] = 4
$data[78,4] = 369
$data[78,5] = 9
$data[78,6] = 0
$data[78,7] = 2
$data[79,0] = 'Uruguay'
$data[79,1] = 350
$data[79,2] = 0
$data[79,3] = 62
$data[79,4] = 286
$data[79,5] = 15
$data[79,6] = 0
$data[79,7] = 2
$data[80,0] = 'Kuwait'
$data[80,1] = 342
$data[80,2] = 25
$data[80,3] = 81
$data[80,4] = 261
$data[80,5] = 15
$data[80,6] = 0
$data[80,7] = 0
$data[81,0] = 'Taiwan'
$data[81,1] = 339
$data[81,2] = 10
$data[81,3] = 50
$data[81,4] = 284
$data[81,5] = 0
$data[81,6] = 0
$data[81,7] = 5
$data[82,0] = 'Republica de Chipre'
$data[82,1] = 320
$data[82,2] = 0
$data[82,3] = 28
$data[82,4] = 283
$data[82,5] = 11
$data[82,6] = 0
$data[82,7] = 9
$data[83,0] = 'Burkina Faso'
$data[83,1] = 288
$data[83,2] = 6
$data[83,3] = 50
$data[83,4] = 222
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 16
$data[84,0] = 'Camerun'
$data[84,1] = 284
$data[84,2] = 51
$data[84,3] = 10
$data[84,4] = 267
$data[84,5] = 0
$data[84,6] = 1
$data[84,7] = 7
$data[85,0] = 'Reunion'
$data[85,1] = 281
$data[85,2] = 0
$data[85,3] = 40
$data[85,4] = 241
$data[85,5] = 3
$data[85,6] = 0
$data[85,7] = 0
$data[86,0] = 'Jordania'
$data[86,1] = 278
$data[86,2] = 0
$data[86,3] = 36
$data[86,4] = 237
$data[86,5] = 5
$data[86,6] = 0
$data[86,7] = 5
$data[87,0] = 'Albania'
$data[87,1] = 277
$data[87,2] = 18
$data[87,3] = 67
$data[87,4] = 194
$data[87,5] = 7
$data[87,6] = 1
$data[87,7] = 16
$data[88,0] = 'Afganistan'
$data[88,1] = 239
$data[88,2] = 2
$data[88,3] = 5
$data[88,4] = 230
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 4
$data[89,0] = 'San Marino'
$data[89,1] = 236
$data[89,2] = 0
$data[89,3] = 13
$data[89,4] = 195
$data[89,5] = 16
$data[89,6] = 0
$data[89,7] = 28
$data[90,0] = 'Oman'
$data[90,1] = 231
$data[90,2] = 21
$data[90,3] = 57
$data[90,4] = 173
$data[90,5] = 3
$data[90,6] = 0
$data[90,7] = 1
$data[91,0] = 'Vietnam'
$data[91,1] = 227
$data[91,2] = 9
$data[91,3] = 75
$data[91,4] = 152
$data[91,5] = 3
$data[91,6] = 0
$data[91,7] = 0
$data[92,0] = 'Honduras'
$data[92,1] = 219
$data[92,2] = 47
$data[92,3] = 3
$data[92,4] = 202
$data[92,5] = 4
$data[92,6] = 4
$data[92,7] = 14
$data[93,0] = 'Cuba'
$data[93,1] = 212
$data[93,2] = 0
$data[93,3] = 12
$data[93,4] = 194
$data[93,5] = 3
$data[93,6] = 0
$data[93,7] = 6
$data[94,0] = 'Malta'
$data[94,1] = 196
$data[94,2] = 8
$data[94,3] = 2
$data[94,4] = 194
$data[94,5] = 2
$data[94,6] = 0
$data[94,7] = 0
$data[95,0] = 'Ghana'
$data[95,1] = 195
$data[95,2] = 0
$data[95,3] = 31
$data[95,4] = 159
$data[95,5] = 1
$data[95,6] = 0
$data[95,7] = 5
$data[96,0] = 'Senegal'
$data[96,1] = 195
$data[96,2] = 5
$data[96,3] = 55
$data[96,4] = 139
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 1
$data[97,0] = 'Costa de Marfil'
$data[97,1] = 190
$data[97,2] = 0
$data[97,3] = 9
$data[97,4] = 180
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 1
$data[98,0] = 'Uzbekistan'
$data[98,1] = 190
$data[98,2] = 9
$data[98,3] = 25
$data[98,4] = 163
$data[98,5] = 8
$data[98,6] = 0
$data[98,7] = 2
$data[99,0] = 'Islas Feroe'
$data[99,1] = 177
$data[99,2] = 4
$data[99,3] = 81
$data[99,4] = 96
$data[99,5] = 1
$data[99,6] = 0
$data[99,7] = 0
$data[100,0] = 'Nigeria'
$data[100,1] = 174
$data[100,2] = 0
$data[100,3] = 9
$data[100,4] = 163
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 2
$data[101,0] = 'Bielorrusia'
$data[101,1] = 163
$data[101,2] = 0
$data[101,3] = 53
$data[101,4] = 108
$data[101,5] = 2
$data[101,6] = 0
$data[101,7] = 2
$data[102,0] = 'Mauricio'
$data[102,1] = 161
$data[102,2] = 0
$data[102,3] = 0
$data[102,4] = 154
$data[102,5] = 1
$data[102,6] = 1
$data[102,7] = 7
$data[103,0] = 'Estado de Palestina'
$data[103,1] = 155
$data[103,2] = 21
$data[103,3] = 18
$data[103,4] = 136
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 1
$data[104,0] = 'Sri Lanka'
$data[104,1] = 150
$data[104,2] = 4
$data[104,3] = 21
$data[104,4] = 126
$data[104,5] = 5
$data[104,6] = 0
$data[104,7] = 3
$data[105,0] = 'Venezuela'
$data[105,1] = 144
$data[105,2] = 0
$data[105,3] = 43
$data[105,4] = 98
$data[105,5] = 6
$data[105,6] = 0
$data[105,7] = 3
$data[106,0] = 'Montenegro'
$data[106,1] = 140
$data[106,2] = 17
$data[106,3] = 0
$data[106,4] = 138
$data[106,5] = 4
$data[106,6] = 0
$data[106,7] = 2
$data[107,0] = 'Martinica'
$data[107,1] = 135
$data[107,2] = 0
$data[107,3] = 27
$data[107,4] = 105
$data[107,5] = 16
$data[107,6] = 0
$data[107,7] = 3
$data[108,0] = 'Brunei'
$data[108,1] = 133
$data[108,2] = 2
$data[108,3] = 56
$data[108,4] = 76
$data[108,5] = 3
$data[108,6] = 0
$data[108,7] = 1
$data[109,0] = 'Georgia'
$data[109,1] = 130
$data[109,2] = 13
$data[109,3] = 26
$data[109,4] = 104
$data[109,5] = 6
$data[109,6] = 0
$data[109,7] = 0
$data[110,0] = 'Guadalupe'
$data[110,1] = 125
$data[110,2] = 0
$data[110,3] = 24
$data[110,4] = 95
$data[110,5] = 14
$data[110,6] = 0
$data[110,7] = 6
$data[111,0] = 'Bolivia'
$data[111,1] = 123
$data[111,2] = 8
$data[111,3] = 1
$data[111,4] = 115
$data[111,5] = 3
$data[111,6] = 0
$data[111,7] = 7
$data[112,0] = 'Consejo Danes para los Refugiados'
$data[112,1] = 123
$data[112,2] = 14
$data[112,3] = 3
$data[112,4] = 109
$data[112,5] = 0
$data[112,6] = 2
$data[112,7] = 11
$data[113,0] = 'Kirguistan'
$data[113,1] = 116
$data[113,2] = 5
$data[113,3] = 5
$data[113,4] = 111
$data[113,5] = 5
$data[113,6] = 0
$data[113,7] = 0
$data[114,0] = 'Mayotte'
$data[114,1] = 116
$data[114,2] = 15
$data[114,3] = 10
$data[114,4] = 105
$data[114,5] = 3
$data[114,6] = 0
$data[114,7] = 1
$data[115,0] = 'Kenia'
$data[115,1] = 110
$data[115,2] = 29
$data[115,3] = 4
$data[115,4] = 103
$data[115,5] = 2
$data[115,6] = 2
$data[115,7] = 3
$data[116,0] = 'Camboya'
$data[116,1] = 110
$data[116,2] = 1
$data[116,3] = 34
$data[116,4] = 76
$data[116,5] = 1
$data[116,6] = 0
$data[116,7] = 0
$data[117,0] = 'Trinidad yTobago'
$data[117,1] = 90
$data[117,2] = 0
$data[117,3] = 1
$data[117,4] = 84
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 5
$data[118,0] = 'Ruanda'
$data[118,1] = 82
$data[118,2] = 0
$data[118,3] = 0
$data[118,4] = 82
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 0
$data[119,0] = 'Gibraltar'
$data[119,1] = 81
$data[119,2] = 0
$data[119,3] = 34
$data[119,4] = 47
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = 'Paraguay'
$data[120,1] = 77
$data[120,2] = 8
$data[120,3] = 2
$data[120,4] = 72
$data[120,5] = 4
$data[120,6] = 0
$data[120,7] = 3
$data[121,0] = 'Isla de Man'
$data[121,1] = 75
$data[121,2] = 7
$data[121,3] = 0
$data[121,4] = 74
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 1
$data[122,0] = 'Niger'
$data[122,1] = 74
$data[122,2] = 0
$data[122,3] = 0
$data[122,4] = 69
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 5
$data[123,0] = 'Liechtenstein'
$data[123,1] = 72
$data[123,2] = 0
$data[123,3] = 0
$data[123,4] = 72
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = 'Madagascar'
$data[124,1] = 59
$data[124,2] = 2
$data[124,3] = 0
$data[124,4] = 59
$data[124,5] = 6
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = 'Banglades'
$data[125,1] = 56
$data[125,2] = 2
$data[125,3] = 25
$data[125,4] = 25
$data[125,5] = 1
$data[125,6] = 0
$data[125,7] = 6
$data[126,0] = 'Aruba'
$data[126,1] = 55
$data[126,2] = 0
$data[126,3] = 1
$data[126,4] = 54
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 0
$data[127,0] = 'Monaco'
$data[127,1] = 55
$data[127,2] = 0
$data[127,3] = 2
$data[127,4] = 52
$data[127,5] = 2
$data[127,6] = 0
$data[127,7] = 1
$data[128,0] = 'Guinea'
$data[128,1] = 52
$data[128,2] = 22
$data[128,3] = 0
$data[128,4] = 52
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = 'Guayana Francesa'
$data[129,1] = 51
$data[129,2] = 0
$data[129,3] = 15
$data[129,4] = 36
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = 'Guatemala'
$data[130,1] = 46
$data[130,2] = 7
$data[130,3] = 12
$data[130,4] = 33
$data[130,5] = 1
$data[130,6] = 0
$data[130,7] = 1
$data[131,0] = 'Barbados'
$data[131,1] = 45
$data[131,2] = 0
$data[131,3] = 0
$data[131,4] = 45
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 0
$data[132,0] = 'Uganda'
$data[132,1] = 44
$data[132,2] = 0
$data[132,3] = 0
$data[132,4] = 44
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = 'Jamaica'
$data[133,1] = 44
$data[133,2] = 0
$data[133,3] = 2
$data[133,4] = 39
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 3
$data[134,0] = 'El Salvador'
$data[134,1] = 41
$data[134,2] = 8
$data[134,3] = 0
$data[134,4] = 39
$data[134,5] = 4
$data[134,6] = 0
$data[134,7] = 2
$data[135,0] = 'Macao'
$data[135,1] = 41
$data[135,2] = 0
$data[135,3] = 10
$data[135,4] = 31
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = 'Republica de Yibuti'
$data[136,1] = 40
$data[136,2] = 7
$data[136,3] = 0
$data[136,4] = 40
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = 'Zambia'
$data[137,1] = 39
$data[137,2] = 3
$data[137,3] = 0
$data[137,4] = 38
$data[137,5] = 0
$data[137,6] = 1
$data[137,7] = 1
$data[138,0] = 'Puerto Rico'
$data[138,1] = 39
$data[138,2] = 0
$data[138,3] = 1
$data[138,4] = 36
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 2
$data[139,0] = 'Polinesia Francesa'
$data[139,1] = 37
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 37
$data[139,5] = 1
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = 'Togo'
$data[140,1] = 36
$data[140,2] = 0
$data[140,3] = 10
$data[140,4] = 24
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 2
$data[141,0] = 'Guam'
$data[141,1] = 32
$data[141,2] = 0
$data[141,3] = 0
$data[141,4] = 31
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 1
$data[142,0] = 'Bermudas'
$data[142,1] = 32
$data[142,2] = 0
$data[142,3] = 10
$data[142,4] = 22
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = 'Mali'
$data[143,1] = 31
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 28
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 3
$data[144,0] = 'Etiopia'
$data[144,1] = 29
$data[144,2] = 0
$data[144,3] = 3
$data[144,4] = 26
$data[144,5] = 2
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = 'Islas Caimanes'
$data[145,1] = 22
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 21
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 1
$data[146,0] = 'Congo'
$data[146,1] = 22
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 20
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 2
$data[147,0] = 'San Martin (Parte Francesa)'
$data[147,1] = 22
$data[147,2] = 0
$data[147,3] = 2
$data[147,4] = 19
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 1
$data[148,0] = 'Bahamas'
$data[148,1] = 21
$data[148,2] = 0
$data[148,3] = 1
$data[148,4] = 19
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 1
$data[149,0] = 'Tanzania'
$data[149,1] = 20
$data[149,2] = 0
$data[149,3] = 2
$data[149,4] = 17
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 1
$data[150,0] = 'Guyana'
$data[150,1] = 19
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 15
$data[150,5] = 0
$data[150,6] = 1
$data[150,7] = 4
$data[151,0] = 'Maldivas'
$data[151,1] = 19
$data[151,2] = 0
$data[151,3] = 13
$data[151,4] = 6
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 0
$data[152,0] = 'Eritrea'
$data[152,1] = 18
$data[152,2] = 3
$data[152,3] = 0
$data[152,4] = 18
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = 'Gabon'
$data[153,1] = 18
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 17
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 1
$data[154,0] = 'Nueva Caledonia'
$data[154,1] = 18
$data[154,2] = 2
$data[154,3] = 1
$data[154,4] = 17
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = 'Islas Virgenes de los Estados Unidos'
$data[155,1] = 17
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 17
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = 'Haiti'
$data[156,1] = 16
$data[156,2] = 0
$data[156,3] = 1
$data[156,4] = 15
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = 'Birmania'
$data[157,1] = 16
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 15
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 1
$data[158,0] = 'San Martin (Parte Holandesa)'
$data[158,1] = 16
$data[158,2] = 0
$data[158,3] = 6
$data[158,4] = 9
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 1
$data[159,0] = 'Guinea Ecuatorial'
$data[159,1] = 15
$data[159,2] = 0
$data[159,3] = 1
$data[159,4] = 14
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = 'Mongolia'
$data[160,1] = 14
$data[160,2] = 0
$data[160,3] = 2
$data[160,4] = 12
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = 'Namibia'
$data[161,1] = 14
$data[161,2] = 0
$data[161,3] = 2
$data[161,4] = 12
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = 'Santa Lucia'
$data[162,1] = 13
$data[162,2] = 0
$data[162,3] = 1
$data[162,4] = 12
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = 'Benin'
$data[163,1] = 13
$data[163,2] = 0
$data[163,3] = 1
$data[163,4] = 12
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = 'Dominica'
$data[164,1] = 12
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 12
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = 'Curazao'
$data[165,1] = 11
$data[165,2] = 0
$data[165,3] = 3
$data[165,4] = 7
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 1
$data[166,0] = 'Surinam'
$data[166,1] = 10
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 10
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = 'Libia'
$data[167,1] = 10
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 10
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = 'Granada'
$data[168,1] = 10
$data[168,2] = 1
$data[168,3] = 0
$data[168,4] = 10
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = 'Laos'
$data[169,1] = 10
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 10
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = 'Mozambique'
$data[170,1] = 10
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 10
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = 'Seychelles'
$data[171,1] = 10
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 10
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Groenlandia'
$data[172,1] = 10
$data[172,2] = 0
$data[172,3] = 2
$data[172,4] = 8
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = 'Siria'
$data[173,1] = 10
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 8
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 2
$data[174,0] = 'Suazilandia'
$data[174,1] = 9
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 9
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = 'Guinea-Bisau'
$data[175,1] = 9
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 9
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = 'Montserrat'
$data[176,1] = 9
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 7
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 2
$data[177,0] = 'San Cristobal y Nieves'
$data[177,1] = 8
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 8
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = 'Zimbabue'
$data[178,1] = 8
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 7
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 1
$data[179,0] = 'Angola'
$data[179,1] = 8
$data[179,2] = 0
$data[179,3] = 1
$data[179,4] = 5
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 2
$data[180,0] = 'Republica del Chad'
$data[180,1] = 7
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 7
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = 'Antigua y Barbuda'
$data[181,1] = 7
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 7
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = 'Fiyi'
$data[182,1] = 7
$data[182,2] = 2
$data[182,3] = 0
$data[182,4] = 7
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Sudan'
$data[183,1] = 7
$data[183,2] = 0
$data[183,3] = 2
$data[183,4] = 3
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 2
$data[184,0] = 'Santa Sede'
$data[184,1] = 6
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 6
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Liberia'
$data[185,1] = 6
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 6
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Islas Turcas y Caicos'
$data[186,1] = 6
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 6
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = 'Nepal'
$data[187,1] = 6
$data[187,2] = 1
$data[187,3] = 1
$data[187,4] = 5
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = 'Cabo Verde'
$data[188,1] = 6
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 5
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 1
$data[189,0] = 'San Bartolome'
$data[189,1] = 6
$data[189,2] = 0
$data[189,3] = 1
$data[189,4] = 5
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = 'Mauritania'
$data[190,1] = 6
$data[190,2] = 0
$data[190,3] = 2
$data[190,4] = 3
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 1
$data[191,0] = 'Nicaragua'
$data[191,1] = 5
$data[191,2] = 0
$data[191,3] = 0
$data[191,4] = 4
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 1
$data[192,0] = 'Somalia'
$data[192,1] = 5
$data[192,2] = 0
$data[192,3] = 1
$data[192,4] = 4
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = 'Butan'
$data[193,1] = 5
$data[193,2] = 1
$data[193,3] = 1
$data[193,4] = 4
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = 'Botsuana'
$data[194,1] = 4
$data[194,2] = 0
$data[194,3] = 0
$data[194,4] = 3
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 1
$data[195,0] = 'Gambia'
$data[195,1] = 4
$data[195,2] = 0
$data[195,3] = 2
$data[195,4] = 1
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 1
$data[196,0] = 'Belice'
$data[196,1] = 3
$data[196,2] = 0
$data[196,3] = 0
$data[196,4] = 3
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = 'Republica de Africa Central'
$data[197,1] = 3
$data[197,2] = 0
$data[197,3] = 0
$data[197,4] = 3
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = 'Islas Virgenes Britanicas'
$data[198,1] = 3
$data[198,2] = 0
$data[198,3] = 0
$data[198,4] = 3
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = 'Bonaire, San Eustaquio y Saba'
$data[199,1] = 2
$data[199,2] = 0
$data[199,3] = 0
$data[199,4] = 2
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = 'Anguila'
$data[200,1] = 2
$data[200,2] = 0
$data[200,3] = 0
$data[200,4] = 2
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = 'Sierra Leona'
$data[201,1] = 2
$data[201,2] = 0
$data[201,3] = 0
$data[201,4] = 2
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = 'Burundi'
$data[202,1] = 2
$data[202,2] = 0
$data[202,3] = 0
$data[202,4] = 2
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = 'San Vicente y las Granadinas'
$data[203,1] = 2
$data[203,2] = 0
$data[203,3] = 1
$data[203,4] = 1
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = 'Papua Nueva Guinea'
$data[204,1] = 1
$data[204,2] = 0
$data[204,3] = 0
$data[204,4] = 1
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = 'Timor Oriental'
$data[205,1] = 1
$data[205,2] = 0
$data[205,3] = 0
$data[205,4] = 1
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0

$ws.Range("A4:H209").Value = $data
